$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 325.4
$ws.Range("I33").Value = 299.2143
$ws.Range("J33").Value = 386.5
$ws.Range("K33").Value = 299.2143
$ws.Range("L33").Value = 386.5
$ws.Range("M33").Value = -70.21429999999998
$ws.Range("N33").Value = -844.5
$ws.Range("H112").Value = 27779412
$ws.Range("I112").Value = 333333760
$ws.Range("J112").Value = 3004736.5
$ws.Range("K112").Value = 1000001280
$ws.Range("L112").Value = 9014209.5
$ws.Range("M112").Value = -1000000172
$ws.Range("N112").Value = -9016425.5
$ws.Range("H116").Value = 1877.3334
$ws.Range("J116").Value = 2151.5
$ws.Range("L116").Value = 2151.5
$ws.Range("N116").Value = -9035.5
$ws.Range("H137").Value = 3573494.2
$ws.Range("I137").Value = 4168420.2
$ws.Range("J137").Value = 3937.5
$ws.Range("K137").Value = 12505260.6
$ws.Range("L137").Value = 11812.5
$ws.Range("M137").Value = -12502710.6
$ws.Range("N137").Value = -16912.5
$ws.Range("H138").Value = 2690977.5
$ws.Range("I138").Value = 1474.2916
$ws.Range("J138").Value = 4389611
$ws.Range("K138").Value = 4422.8748
$ws.Range("L138").Value = 13168833
$ws.Range("M138").Value = 717.1252000000004
$ws.Range("N138").Value = -13179113
$ws.Range("H141").Value = 1237.875
$ws.Range("I141").Value = 617.1667
$ws.Range("J141").Value = 3100
$ws.Range("K141").Value = 1851.5001
$ws.Range("L141").Value = 9300
$ws.Range("M141").Value = 3328.4999
$ws.Range("N141").Value = -19660

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 13158186
$ws.Range("I5").Value = 17543914
$ws.Range("K5").Value = 17543914
$ws.Range("M5").Value = -17543802
$ws.Range("H45").Value = 1957.5
$ws.Range("I45").Value = 2840
$ws.Range("J45").Value = 1327.1428
$ws.Range("K45").Value = 2840
$ws.Range("L45").Value = 1327.1428
$ws.Range("M45").Value = -2463
$ws.Range("N45").Value = -2081.1428
$ws.Range("H109").Value = 43000
$ws.Range("J109").Value = 43000
$ws.Range("L109").Value = 43000
$ws.Range("N109").Value = -45774
$ws.Range("H122").Value = 3833201.8
$ws.Range("I122").Value = 1810.3636
$ws.Range("K122").Value = 5431.0908
$ws.Range("M122").Value = -2981.0908
$ws.Range("H132").Value = 36194.895
$ws.Range("I132").Value = 25151.049
$ws.Range("J132").Value = 64494.75
$ws.Range("K132").Value = 75453.147
$ws.Range("L132").Value = 193484.25
$ws.Range("M132").Value = -72923.147
$ws.Range("N132").Value = -198544.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 13158186
$ws.Range("I4").Value = 17543914
$ws.Range("K4").Value = 17543914
$ws.Range("M4").Value = -17543799
$ws.Range("H22").Value = 176.6
$ws.Range("I22").Value = 196.07692
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 196.07692
$ws.Range("L22").Value = 50
$ws.Range("M22").Value = -23.07692
$ws.Range("N22").Value = -396
$ws.Range("H86").Value = 19640.525
$ws.Range("I86").Value = 16572.25
$ws.Range("K86").Value = 16572.25
$ws.Range("M86").Value = -15449.25
$ws.Range("H89").Value = 19640.525
$ws.Range("I89").Value = 16572.25
$ws.Range("K89").Value = 82861.25
$ws.Range("M89").Value = -77245.25
$ws.Range("H134").Value = 2216.356
$ws.Range("I134").Value = 1743.1041
$ws.Range("J134").Value = 4281.4546
$ws.Range("K134").Value = 5229.3123
$ws.Range("L134").Value = 12844.3638
$ws.Range("M134").Value = -2694.3123
$ws.Range("N134").Value = -17914.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1730.4789
$ws.Range("I31").Value = 998.93335
$ws.Range("J31").Value = 2996.6155
$ws.Range("K31").Value = 998.93335
$ws.Range("L31").Value = 2996.6155
$ws.Range("M31").Value = -703.93335
$ws.Range("N31").Value = -3586.6155
$ws.Range("H34").Value = 1730.4789
$ws.Range("I34").Value = 998.93335
$ws.Range("J34").Value = 2996.6155
$ws.Range("K34").Value = 998.93335
$ws.Range("L34").Value = 2996.6155
$ws.Range("M34").Value = -796.93335
$ws.Range("N34").Value = -3400.6155
$ws.Range("H58").Value = 22728858
$ws.Range("I58").Value = 25001342
$ws.Range("J58").Value = 4000.25
$ws.Range("K58").Value = 25001342
$ws.Range("L58").Value = 4000.25
$ws.Range("M58").Value = -25001139
$ws.Range("N58").Value = -4406.25
$ws.Range("H107").Value = 363.71875
$ws.Range("I107").Value = 350.42856
$ws.Range("J107").Value = 389.0909
$ws.Range("K107").Value = 350.42856
$ws.Range("L107").Value = 389.0909
$ws.Range("M107").Value = 1569.57144
$ws.Range("N107").Value = -4229.0909
$ws.Range("H134").Value = 24790.877
$ws.Range("I134").Value = 1881
$ws.Range("J134").Value = 114139.4
$ws.Range("K134").Value = 5643
$ws.Range("L134").Value = 342418.2
$ws.Range("M134").Value = -3108
$ws.Range("N134").Value = -347488.2
$ws.Range("H136").Value = 22728858
$ws.Range("I136").Value = 25001342
$ws.Range("J136").Value = 4000.25
$ws.Range("K136").Value = 75004026
$ws.Range("L136").Value = 12000.75
$ws.Range("M136").Value = -75001476
$ws.Range("N136").Value = -17100.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1008.1205
$ws.Range("I131").Value = 663.1667
$ws.Range("J131").Value = 1035
$ws.Range("K131").Value = 1989.5001
$ws.Range("L131").Value = 3105
$ws.Range("M131").Value = 3050.4999
$ws.Range("N131").Value = -13185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 916.75
$ws.Range("I102").Value = 844.55554
$ws.Range("K102").Value = 844.55554
$ws.Range("M102").Value = 777.44446
$ws.Range("H132").Value = 64256.906
$ws.Range("I132").Value = 44942.305
$ws.Range("J132").Value = 113616.445
$ws.Range("K132").Value = 134826.915
$ws.Range("L132").Value = 340849.335
$ws.Range("M132").Value = -132296.915
$ws.Range("N132").Value = -345909.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 71999.484
$ws.Range("I136").Value = 39845.58
$ws.Range("J136").Value = 350666.66
$ws.Range("K136").Value = 119536.74
$ws.Range("L136").Value = 1051999.98
$ws.Range("M136").Value = -116986.74
$ws.Range("N136").Value = -1057099.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 13467.667
$ws.Range("I96").Value = 10951.5
$ws.Range("J96").Value = 18500
$ws.Range("K96").Value = 10951.5
$ws.Range("L96").Value = 18500
$ws.Range("M96").Value = -9578.5
$ws.Range("N96").Value = -21246
$ws.Range("H107").Value = 319.5
$ws.Range("I107").Value = 287
$ws.Range("J107").Value = 352
$ws.Range("K107").Value = 861
$ws.Range("L107").Value = 1056
$ws.Range("M107").Value = 1059
$ws.Range("N107").Value = -4896
$ws.Range("H136").Value = 39027.527
$ws.Range("I136").Value = 21376.715
$ws.Range("K136").Value = 64130.145
$ws.Range("M136").Value = -61580.145
